$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 4): lower-case "polXX" -> capitalised "PolXX"
$ws.Range("A4").Value = "Pol0"
$ws.Range("B4").Value = "Pol45"
$ws.Range("C4").Value = "Pol90"
$ws.Range("D4").Value = "Pol135"

# Data rows (5-10): fix the path by inserting the missing
# "SampleImageSetByExcelFileFinder" path segment before "FourCamera"
$oldSegment = "/excel/FourCamera/"
$newSegment = "/excel/SampleImageSetByExcelFileFinder/FourCamera/"

for ($row = 5; $row -le 10; $row++) {
    for ($col = 1; $col -le 4; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $oldValue = $cell.Value()
        $cell.Value = $oldValue.Replace($oldSegment, $newSegment)
    }
}

# Move the active selection from D11 to D5
$ws.Range("D5").Select()
